# "removed matched data from excel report"
# Replace the header row labels (A1:D1) with the new report column
# headers. E1 ("df") is left untouched since it did not change.
#
# Values are written in B1, C1, D1, A1 order (rather than left-to-right)
# so the workbook's shared-string table is rebuilt in the same order as
# the target file (Risk, Cell, Mask, Trade ID appended after the
# untouched df/sef/hjhj strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Cell"
$ws.Range("D1").Value = "Mask"
$ws.Range("A1").Value = "Trade ID"

# Reset the active selection back to A1 (closest achievable match to the
# target file, which no longer carries an explicit B3 selection).
$ws.Range("A1").Select()
